$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação date update (keep stored as text, matching original "01/01/2018" string)
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2021"
$ws.Range("C8").Value = "01/01/2021"

# Programa (Portuguese) text update - spacing normalized, last sentence removed
$ws.Range("B16").Value = "1. Teoria dos Sistemas. 1.1. Conceitos Básicos sobre Teoria dos Sistemas;1.2. Abordagem Sistêmica; 1.3. Classificação dos Sistemas; 1.4. Ciclo de Vida dos Sistemas; 2. Processo de Modelagem de Sistemas. 2.1. Modelo Formal; 2.2. Modelo Computacional; 2.3. Teoria dos Modelos; 3. Modelagem para Simulação 3.1. Modelagem Discreta; 3.2. Modelagem Orientada a Eventos; 4. Linguagens de Simulação. 5. Aspectos Estatísticos da Simulação de Sistemas. 5.1. Geradores de Números Aleatórios; 5.2. Geração de Variáveis Aleatórias; 5.3. Inferência Estatística; 5.4. Problemas Estatísticos Relacionados com Simulação; 6. Validação de Modelos."
$ws.Range("C16").Value = "1. Teoria dos Sistemas. 1.1. Conceitos Básicos sobre Teoria dos Sistemas;1.2. Abordagem Sistêmica; 1.3. Classificação dos Sistemas; 1.4. Ciclo de Vida dos Sistemas; 2. Processo de Modelagem de Sistemas. 2.1. Modelo Formal; 2.2. Modelo Computacional; 2.3. Teoria dos Modelos; 3. Modelagem para Simulação 3.1. Modelagem Discreta; 3.2. Modelagem Orientada a Eventos; 4. Linguagens de Simulação. 5. Aspectos Estatísticos da Simulação de Sistemas. 5.1. Geradores de Números Aleatórios; 5.2. Geração de Variáveis Aleatórias; 5.3. Inferência Estatística; 5.4. Problemas Estatísticos Relacionados com Simulação; 6. Validação de Modelos."

# Syllabus (English) text - previously wrongly duplicated the short syllabus text; now filled with full English syllabus
$ws.Range("B17").Value = "1. Systems Theory.1.1. Basic Concepts on Systems Theory;1.2. Systemic Approach;1.3. Classification of Systems;1.4. Systems Life Cycle;2. System Modeling Process.2.1. Formal Model;2.2. Computational Model;2.3. Theory of Models;3. Modeling for Simulation3.1. Discrete Modeling;3.2. Event-Driven Modeling;4. Simulation Languages.5. Statistical Aspects of Systems Simulation.5.1. Random Number Generators;5.2. Generation of Random Variables;5.3. Statistical inference;5.4. Statistical Problems Related to Simulation;6. Validation of Models."
$ws.Range("C17").Value = "1. Systems Theory.1.1. Basic Concepts on Systems Theory;1.2. Systemic Approach;1.3. Classification of Systems;1.4. Systems Life Cycle;2. System Modeling Process.2.1. Formal Model;2.2. Computational Model;2.3. Theory of Models;3. Modeling for Simulation3.1. Discrete Modeling;3.2. Event-Driven Modeling;4. Simulation Languages.5. Statistical Aspects of Systems Simulation.5.1. Random Number Generators;5.2. Generation of Random Variables;5.3. Statistical inference;5.4. Statistical Problems Related to Simulation;6. Validation of Models."

# Método (Method) text update
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Critério (Criterion) text update
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# Norma de recuperação (Recovery norm) text update
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."

# Bibliografia (Bibliography) text update
$ws.Range("B22").Value = "1. BANKS, J., CARSON, J.S., “Discrete-Event System Simulation.”, Printice Hall, 2000.2. CHWIF, L., MEDINA, A.C. e col. “Introdução ao Simul8: um guia prático”, Livro Eletrônico, 1ª ed., 2015.3. Fishman, G.S. “Discrete-Event Simulation: Modeling, Programming, and Analysis”, Springer-Verlag, 2001."
$ws.Range("C22").Value = "1. BANKS, J., CARSON, J.S., “Discrete-Event System Simulation.”, Printice Hall, 2000.2. CHWIF, L., MEDINA, A.C. e col. “Introdução ao Simul8: um guia prático”, Livro Eletrônico, 1ª ed., 2015.3. Fishman, G.S. “Discrete-Event Simulation: Modeling, Programming, and Analysis”, Springer-Verlag, 2001."
